$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price values that must remain text
$textCells = @("D5", "D8", "D10", "D11", "D15", "D17", "D20", "D22", "D23", "D25", "D26", "D28", "D29", "D32", "D41", "D46", "D48")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply cell value updates from the diff
$ws.Range("D2").Value = "26.009.53"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.636.40"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "214.58"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "0.251"
$ws.Range("E8").Value = "  -1.91%  "
$ws.Range("E9").Value = "  -2.95%  "
$ws.Range("D10").Value = "18.27"
$ws.Range("E10").Value = "  -6.96%  "
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "1.864.25"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("D14").Value = "1.627.88"
$ws.Range("E14").Value = "  -3.10%  "
$ws.Range("D15").Value = "0.527"
$ws.Range("E15").Value = "  -3.15%  "
$ws.Range("D16").Value = "26.001.80"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "61.60"
$ws.Range("E17").Value = "  -2.69%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  -3.26%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "191.08"
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("E21").Value = "  -2.35%  "
$ws.Range("D22").Value = "9.69"
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("D23").Value = "6.08"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("D25").Value = "143.83"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").Value = "1.78"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "6.79"
$ws.Range("E28").Value = "  -1.57%  "
$ws.Range("D29").Value = "15.18"
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("E30").Value = "  -1.35%  "
$ws.Range("E31").Value = "  -2.91%  "
$ws.Range("D32").Value = "3.14"
$ws.Range("E32").Value = "  -3.75%  "
$ws.Range("E33").Value = "  -4.78%  "
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("E35").Value = "  -2.84%  "
$ws.Range("D36").Value = "1.135.64"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  -4.50%  "
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("E39").Value = "  -4.04%  "
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("D41").Value = "98.38"
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("E42").Value = "  -2.63%  "
$ws.Range("D43").Value = "1.774.86"
$ws.Range("E44").Value = "  -4.81%  "
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").Value = "55.03"
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").Value = "1.48"
$ws.Range("E48").Value = "  +2.10%  "
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("E51").Value = "  +0.10%  "
